$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had only the "sCs" sending-cluster block (rows 2-6).
# Following Dr Hou's advice, a second sending cluster "FAPs" is added with
# the same 5 target clusters (ECs, FAPs, M1, M2, sCs), and the numeric
# values for all rows are refreshed with the updated computation.
#
# Resulting layout (rows 2-11), columns A-D are text, E-T are numeric:
#   Row 2 : FAPs | Wnt7b | Fzd5 | ECs
#   Row 3 : FAPs | Wnt7b | Fzd5 | FAPs
#   Row 4 : FAPs | Wnt7b | Fzd5 | M1
#   Row 5 : FAPs | Wnt7b | Fzd5 | M2
#   Row 6 : FAPs | Wnt7b | Fzd5 | sCs
#   Row 7 : sCs  | Wnt7b | Fzd5 | ECs
#   Row 8 : sCs  | Wnt7b | Fzd5 | FAPs
#   Row 9 : sCs  | Wnt7b | Fzd5 | M1
#   Row 10: sCs  | Wnt7b | Fzd5 | M2
#   Row 11: sCs  | Wnt7b | Fzd5 | sCs

# Row 2: FAPs/Wnt7b/Fzd5/ECs
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Wnt7b"
$ws.Cells.Item(2,3).Value = "Fzd5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.072919
$ws.Cells.Item(2,8).Value = 0.218757
$ws.Cells.Item(2,9).Value = 0.1016383815134179
$ws.Cells.Item(2,10).Value = 0.1016383815134179
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 4.048438999999999
$ws.Cells.Item(2,14).Value = 12.145317
$ws.Cells.Item(2,15).Value = 0.1703267688113503
$ws.Cells.Item(2,16).Value = 0.1703267688113503
$ws.Cells.Item(2,17).Value = 0.2952081234409999
$ws.Cells.Item(2,18).Value = 2.656873110969
$ws.Cells.Item(2,19).Value = 0.01731173711039575
$ws.Cells.Item(2,20).Value = 0.01731173711039576

# Row 3: FAPs/Wnt7b/Fzd5/FAPs
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Wnt7b"
$ws.Cells.Item(3,3).Value = "Fzd5"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.072919
$ws.Cells.Item(3,8).Value = 0.218757
$ws.Cells.Item(3,9).Value = 0.1016383815134179
$ws.Cells.Item(3,10).Value = 0.1016383815134179
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 11.27122833333333
$ws.Cells.Item(3,14).Value = 33.813685
$ws.Cells.Item(3,15).Value = 0.4742054659960562
$ws.Cells.Item(3,16).Value = 0.4742054659960562
$ws.Cells.Item(3,17).Value = 0.8218866988383333
$ws.Cells.Item(3,18).Value = 7.396980289545001
$ws.Cells.Item(3,19).Value = 0.04819747606865529
$ws.Cells.Item(3,20).Value = 0.04819747606865529

# Row 4: FAPs/Wnt7b/Fzd5/M1
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Wnt7b"
$ws.Cells.Item(4,3).Value = "Fzd5"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.072919
$ws.Cells.Item(4,8).Value = 0.218757
$ws.Cells.Item(4,9).Value = 0.1016383815134179
$ws.Cells.Item(4,10).Value = 0.1016383815134179
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.219226333333333
$ws.Cells.Item(4,14).Value = 9.657679
$ws.Cells.Item(4,15).Value = 0.1354399607920677
$ws.Cells.Item(4,16).Value = 0.1354399607920676
$ws.Cells.Item(4,17).Value = 0.2347427650003333
$ws.Cells.Item(4,18).Value = 2.112684885003
$ws.Cells.Item(4,19).Value = 0.01376589840714654
$ws.Cells.Item(4,20).Value = 0.01376589840714654

# Row 5: FAPs/Wnt7b/Fzd5/M2
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Wnt7b"
$ws.Cells.Item(5,3).Value = "Fzd5"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.072919
$ws.Cells.Item(5,8).Value = 0.218757
$ws.Cells.Item(5,9).Value = 0.1016383815134179
$ws.Cells.Item(5,10).Value = 0.1016383815134179
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 4.235549333333334
$ws.Cells.Item(5,14).Value = 12.706648
$ws.Cells.Item(5,15).Value = 0.1781989137264352
$ws.Cells.Item(5,16).Value = 0.1781989137264352
$ws.Cells.Item(5,17).Value = 0.3088520218373334
$ws.Cells.Item(5,18).Value = 2.779668196536
$ws.Cells.Item(5,19).Value = 0.01811184917860407
$ws.Cells.Item(5,20).Value = 0.01811184917860407

# Row 6: FAPs/Wnt7b/Fzd5/sCs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Wnt7b"
$ws.Cells.Item(6,3).Value = "Fzd5"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.072919
$ws.Cells.Item(6,8).Value = 0.218757
$ws.Cells.Item(6,9).Value = 0.1016383815134179
$ws.Cells.Item(6,10).Value = 0.1016383815134179
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.9942166666666665
$ws.Cells.Item(6,14).Value = 2.98265
$ws.Cells.Item(6,15).Value = 0.04182889067409059
$ws.Cells.Item(6,16).Value = 0.04182889067409059
$ws.Cells.Item(6,17).Value = 0.07249728511666666
$ws.Cells.Item(6,18).Value = 0.6524755660499999
$ws.Cells.Item(6,19).Value = 0.004251420748616268
$ws.Cells.Item(6,20).Value = 0.004251420748616268

# Row 7: sCs/Wnt7b/Fzd5/ECs
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Wnt7b"
$ws.Cells.Item(7,3).Value = "Fzd5"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.6445166666666666
$ws.Cells.Item(7,8).Value = 1.93355
$ws.Cells.Item(7,9).Value = 0.8983616184865821
$ws.Cells.Item(7,10).Value = 0.898361618486582
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 4.048438999999999
$ws.Cells.Item(7,14).Value = 12.145317
$ws.Cells.Item(7,15).Value = 0.1703267688113503
$ws.Cells.Item(7,16).Value = 0.1703267688113503
$ws.Cells.Item(7,17).Value = 2.609286409483333
$ws.Cells.Item(7,18).Value = 23.48357768535
$ws.Cells.Item(7,19).Value = 0.1530150317009545
$ws.Cells.Item(7,20).Value = 0.1530150317009545

# Row 8: sCs/Wnt7b/Fzd5/FAPs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Wnt7b"
$ws.Cells.Item(8,3).Value = "Fzd5"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.6445166666666666
$ws.Cells.Item(8,8).Value = 1.93355
$ws.Cells.Item(8,9).Value = 0.8983616184865821
$ws.Cells.Item(8,10).Value = 0.898361618486582
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 11.27122833333333
$ws.Cells.Item(8,14).Value = 33.813685
$ws.Cells.Item(8,15).Value = 0.4742054659960562
$ws.Cells.Item(8,16).Value = 0.4742054659960562
$ws.Cells.Item(8,17).Value = 7.264494514638888
$ws.Cells.Item(8,18).Value = 65.38045063174999
$ws.Cells.Item(8,19).Value = 0.4260079899274009
$ws.Cells.Item(8,20).Value = 0.4260079899274009

# Row 9: sCs/Wnt7b/Fzd5/M1
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Wnt7b"
$ws.Cells.Item(9,3).Value = "Fzd5"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.6445166666666666
$ws.Cells.Item(9,8).Value = 1.93355
$ws.Cells.Item(9,9).Value = 0.8983616184865821
$ws.Cells.Item(9,10).Value = 0.898361618486582
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.219226333333333
$ws.Cells.Item(9,14).Value = 9.657679
$ws.Cells.Item(9,15).Value = 0.1354399607920677
$ws.Cells.Item(9,16).Value = 0.1354399607920676
$ws.Cells.Item(9,17).Value = 2.074845025605555
$ws.Cells.Item(9,18).Value = 18.67360523045
$ws.Cells.Item(9,19).Value = 0.1216740623849211
$ws.Cells.Item(9,20).Value = 0.1216740623849211

# Row 10: sCs/Wnt7b/Fzd5/M2
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Wnt7b"
$ws.Cells.Item(10,3).Value = "Fzd5"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.6445166666666666
$ws.Cells.Item(10,8).Value = 1.93355
$ws.Cells.Item(10,9).Value = 0.8983616184865821
$ws.Cells.Item(10,10).Value = 0.898361618486582
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 4.235549333333334
$ws.Cells.Item(10,14).Value = 12.706648
$ws.Cells.Item(10,15).Value = 0.1781989137264352
$ws.Cells.Item(10,16).Value = 0.1781989137264352
$ws.Cells.Item(10,17).Value = 2.729882137822222
$ws.Cells.Item(10,18).Value = 24.5689392404
$ws.Cells.Item(10,19).Value = 0.1600870645478312
$ws.Cells.Item(10,20).Value = 0.1600870645478311

# Row 11: sCs/Wnt7b/Fzd5/sCs
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Wnt7b"
$ws.Cells.Item(11,3).Value = "Fzd5"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.6445166666666666
$ws.Cells.Item(11,8).Value = 1.93355
$ws.Cells.Item(11,9).Value = 0.8983616184865821
$ws.Cells.Item(11,10).Value = 0.898361618486582
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.9942166666666665
$ws.Cells.Item(11,14).Value = 2.98265
$ws.Cells.Item(11,15).Value = 0.04182889067409059
$ws.Cells.Item(11,16).Value = 0.04182889067409059
$ws.Cells.Item(11,17).Value = 0.6407892119444443
$ws.Cells.Item(11,18).Value = 5.767102907499999
$ws.Cells.Item(11,19).Value = 0.03757746992547432
$ws.Cells.Item(11,20).Value = 0.03757746992547432
